$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 135.71428
$ws.Range("J4").Value = 85
$ws.Range("L4").Value = 85
$ws.Range("N4").Value = -313

$ws.Range("H19").Value = 2260.889
$ws.Range("J19").Value = 2325
$ws.Range("L19").Value = 2325
$ws.Range("N19").Value = -2675

$ws.Range("H28").Value = 852.5
$ws.Range("I28").Value = 852.5
$ws.Range("K28").Value = 852.5
$ws.Range("M28").Value = -367.5

$ws.Range("H107").Value = 222.33333
$ws.Range("I107").Value = 188.75
$ws.Range("J107").Value = 289.5
$ws.Range("K107").Value = 188.75
$ws.Range("L107").Value = 289.5
$ws.Range("M107").Value = 1731.25
$ws.Range("N107").Value = -4129.5

$ws.Range("H137").Value = 2101.5
$ws.Range("J137").Value = 2263.4285
$ws.Range("L137").Value = 6790.2855
$ws.Range("N137").Value = -11890.2855

$ws.Range("H138").Value = 3535.9688
$ws.Range("I138").Value = 2499.6667
$ws.Range("K138").Value = 7499.000100000001
$ws.Range("M138").Value = -2359.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2273.1667
$ws.Range("I45").Value = 1820
$ws.Range("K45").Value = 1820
$ws.Range("M45").Value = -1443

$ws.Range("H61").Value = 5244.75
$ws.Range("I61").Value = 6988.3335
$ws.Range("J61").Value = 14
$ws.Range("K61").Value = 6988.3335
$ws.Range("L61").Value = 14
$ws.Range("M61").Value = -6776.3335
$ws.Range("N61").Value = -438

$ws.Range("H74").Value = 3021.6
$ws.Range("I74").Value = 3051.7144
$ws.Range("K74").Value = 3051.7144
$ws.Range("M74").Value = -2177.7144

$ws.Range("H77").Value = 3021.6
$ws.Range("I77").Value = 3051.7144
$ws.Range("K77").Value = 15258.572
$ws.Range("M77").Value = -10890.572

$ws.Range("H122").Value = 1424.125
$ws.Range("I122").Value = 1399
$ws.Range("J122").Value = 1499.5
$ws.Range("K122").Value = 4197
$ws.Range("L122").Value = 4498.5
$ws.Range("M122").Value = -1747
$ws.Range("N122").Value = -9398.5

$ws.Range("H132").Value = 1569.619
$ws.Range("I132").Value = 1468.4117
$ws.Range("J132").Value = 1999.75
$ws.Range("K132").Value = 4405.2351
$ws.Range("L132").Value = 5999.25
$ws.Range("M132").Value = -1875.2351
$ws.Range("N132").Value = -11059.25

$ws.Range("H136").Value = 5244.75
$ws.Range("I136").Value = 6988.3335
$ws.Range("J136").Value = 14
$ws.Range("K136").Value = 20965.0005
$ws.Range("L136").Value = 42
$ws.Range("M136").Value = -18415.0005
$ws.Range("N136").Value = -5142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2150.7778
$ws.Range("I94").Value = 2150.7778
$ws.Range("K94").Value = 2150.7778
$ws.Range("M94").Value = -1699.7778

$ws.Range("H133").Value = 40000
$ws.Range("I133").Value = 40000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 40000
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -34940
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1850.6364
$ws.Range("I31").Value = 1386.7561
$ws.Range("J31").Value = 3209.1428
$ws.Range("K31").Value = 1386.7561
$ws.Range("L31").Value = 3209.1428
$ws.Range("M31").Value = -1091.7561
$ws.Range("N31").Value = -3799.1428

$ws.Range("H34").Value = 1850.6364
$ws.Range("I34").Value = 1386.7561
$ws.Range("J34").Value = 3209.1428
$ws.Range("K34").Value = 1386.7561
$ws.Range("L34").Value = 3209.1428
$ws.Range("M34").Value = -1184.7561
$ws.Range("N34").Value = -3613.1428

$ws.Range("H58").Value = 5554.8887
$ws.Range("I58").Value = 5011.8125
$ws.Range("K58").Value = 5011.8125
$ws.Range("M58").Value = -4808.8125

$ws.Range("H134").Value = 74087.14
$ws.Range("I134").Value = 101322.2
$ws.Range("K134").Value = 303966.6
$ws.Range("M134").Value = -301431.6

$ws.Range("H136").Value = 5554.8887
$ws.Range("I136").Value = 5011.8125
$ws.Range("K136").Value = 15035.4375
$ws.Range("M136").Value = -12485.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1322
$ws.Range("I75").Value = 644
$ws.Range("J75").Value = 2000
$ws.Range("K75").Value = 1932
$ws.Range("L75").Value = 6000
$ws.Range("M75").Value = -934
$ws.Range("N75").Value = -7996

$ws.Range("H78").Value = 1322
$ws.Range("I78").Value = 644
$ws.Range("J78").Value = 2000
$ws.Range("K78").Value = 5796
$ws.Range("L78").Value = 18000
$ws.Range("M78").Value = -804
$ws.Range("N78").Value = -27984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1892.8
$ws.Range("I132").Value = 1894.75
$ws.Range("J132").Value = 1885
$ws.Range("K132").Value = 5684.25
$ws.Range("L132").Value = 5655
$ws.Range("M132").Value = -3154.25
$ws.Range("N132").Value = -10715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2038.4
$ws.Range("I7").Value = 1842
$ws.Range("J7").Value = 2496.6667
$ws.Range("K7").Value = 1842
$ws.Range("L7").Value = 2496.6667
$ws.Range("M7").Value = -1730
$ws.Range("N7").Value = -2720.6667

$ws.Range("H22").Value = 849.4286
$ws.Range("J22").Value = 716.6667
$ws.Range("L22").Value = 716.6667
$ws.Range("N22").Value = -1306.6667

$ws.Range("H27").Value = 849.4286
$ws.Range("J27").Value = 716.6667
$ws.Range("L27").Value = 716.6667
$ws.Range("N27").Value = -930.6667

$ws.Range("H61").Value = 125
$ws.Range("I61").Value = 125
$ws.Range("K61").Value = 125
$ws.Range("M61").Value = 77

$ws.Range("H113").Value = 125
$ws.Range("I113").Value = 125
$ws.Range("K113").Value = 125
$ws.Range("M113").Value = 2045

$ws.Range("H122").Value = 3808.923
$ws.Range("I122").Value = 3502.889
$ws.Range("J122").Value = 4497.5
$ws.Range("K122").Value = 10508.667
$ws.Range("L122").Value = 13492.5
$ws.Range("M122").Value = -8058.667000000001
$ws.Range("N122").Value = -18392.5

$ws.Range("H126").Value = 2038.4
$ws.Range("I126").Value = 1842
$ws.Range("J126").Value = 2496.6667
$ws.Range("K126").Value = 5526
$ws.Range("L126").Value = 7490.000100000001
$ws.Range("M126").Value = -3056
$ws.Range("N126").Value = -12430.0001

$ws.Range("H132").Value = 5822.75
$ws.Range("I132").Value = 4610.636
$ws.Range("J132").Value = 7304.222
$ws.Range("K132").Value = 13831.908
$ws.Range("L132").Value = 21912.666
$ws.Range("M132").Value = -11301.908
$ws.Range("N132").Value = -26972.666

$ws.Range("H136").Value = 2132.9
$ws.Range("I136").Value = 2132.9
$ws.Range("K136").Value = 6398.700000000001
$ws.Range("M136").Value = -3848.700000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 872
$ws.Range("I132").Value = 727.38464
$ws.Range("K132").Value = 2182.15392
$ws.Range("M132").Value = 347.8460800000003

$ws.Range("H136").Value = 5980.905
$ws.Range("I136").Value = 6251.8
$ws.Range("J136").Value = 5303.6665
$ws.Range("K136").Value = 18755.4
$ws.Range("L136").Value = 15910.9995
$ws.Range("M136").Value = -16205.4
$ws.Range("N136").Value = -21010.9995
